# The "Key"/"ENGLISH"/"FRENCH" header row is removed so the translation
# rows are indexed from the first row (matching a different key order).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the entire first row shifts every row below it up by one,
# which both drops the header row and collapses the now-duplicate last
# row, leaving the used range as A1:C3.
$ws.Rows.Item(1).Delete()
